$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 258 (shifts Rx01.. down by one row)
$ws.Rows.Item(258).Insert()

# Row 257 (O27) grows taller to accommodate the new row being added right below it
$ws.Rows.Item(257).RowHeight = 61

# Fill in the new "Cause of death" derived variable row
$ws.Rows.Item(258).RowHeight = 16
$ws.Cells.Item(258, 1).Value = "O28"
$ws.Cells.Item(258, 2).Value = "cause_of_death"
$ws.Cells.Item(258, 3).Value = "Outcome"
$ws.Cells.Item(258, 4).Value = "Cause of death"
$ws.Cells.Item(258, 5).Value = "1 = COVID-19; 2 = Cancer; 3 = Both; 88 = Other; 99 = Unknown"

$ws.Cells.Item(258, 5).Select()

# Grow the table (ListObject) to include the newly inserted row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E310"))
